$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added ho chi minh" - a new POD location is appended to the list.
# The sheet previously had two trailing blank rows (A28, A29) below the
# last real entry (A27, "Shuwaikh"); the new value goes into the first
# of those (A28), picking up the same vertical-center formatting used by
# the rest of the location rows, and the now-superfluous trailing blank
# row (A29) is removed.

$ws.Range("A28").Value = "Ho Chi Minh City"
$ws.Range("A28").VerticalAlignment = -4108   # xlVAlignCenter - matches A2:A27

# The header cell (A1, "POD") reverts to the plain/default cell style.
$ws.Range("A1").Style = "Normal"

# Drop the trailing blank row that is no longer needed.
$null = $ws.Rows(29).Delete()

# Leave the selection where the edit happened, as the saved file shows.
$null = $ws.Range("B28").Select()
